$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells remain text, matching the source data format
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '66.065.21'
$ws.Range("E2").Value = '  +3.75%  '
$ws.Range("D3").Value = '3.822.38'
$ws.Range("E3").Value = '  +8.56%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '427.37'
$ws.Range("E5").Value = '  +8.07%  '
$ws.Range("D6").Value = '131.54'
$ws.Range("E6").Value = '  +6.53%  '
$ws.Range("D7").Value = '3.816.92'
$ws.Range("E7").Value = '  +8.70%  '
$ws.Range("D8").Value = '0.613'
$ws.Range("E8").Value = '  +3.65%  '
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").Value = '0.734'
$ws.Range("E10").Value = '  +7.42%  '
$ws.Range("D11").Value = '0.158'
$ws.Range("E11").Value = '  +3.35%  '
$ws.Range("D12").Value = '0.0000338'
$ws.Range("E12").Value = '  -0.24%  '
$ws.Range("D13").Value = '41.57'
$ws.Range("E13").Value = '  +6.07%  '
$ws.Range("D14").Value = '10.42'
$ws.Range("E14").Value = '  +12.61%  '
$ws.Range("D15").Value = '4.433.10'
$ws.Range("E15").Value = '  +8.95%  '
$ws.Range("D16").Value = '15.47'
$ws.Range("E16").Value = '  +21.85%  '
$ws.Range("D17").Value = '3.865.72'
$ws.Range("E17").Value = '  +10.11%  '
$ws.Range("D18").Value = '0.138'
$ws.Range("E18").Value = '  +1.44%  '
$ws.Range("D19").Value = '20.05'
$ws.Range("E19").Value = '  +6.48%  '
$ws.Range("D20").Value = '1.11'
$ws.Range("E20").Value = '  +8.10%  '
$ws.Range("D21").Value = '66.238.38'
$ws.Range("E21").Value = '  +4.07%  '
$ws.Range("D22").Value = '417.06'
$ws.Range("E22").Value = '  +4.87%  '
$ws.Range("D23").Value = '15.15'
$ws.Range("E23").Value = '  +9.11%  '
$ws.Range("D24").Value = '85.39'
$ws.Range("E24").Value = '  +4.72%  '
$ws.Range("D25").Value = '3.13'
$ws.Range("E25").Value = '  +9.00%  '
$ws.Range("D26").Value = '37.34'
$ws.Range("E26").Value = '  +10.41%  '
$ws.Range("D27").Value = '10.12'
$ws.Range("E27").Value = '  +14.55%  '
$ws.Range("D28").Value = '3.31'
$ws.Range("E28").Value = '  +10.33%  '
$ws.Range("D29").Value = '9.49'
$ws.Range("E29").Value = '  +38.40%  '
$ws.Range("D30").Value = '5.38'
$ws.Range("E30").Value = '  +2.68%  '
$ws.Range("D31").Value = '14.13'
$ws.Range("E31").Value = '  +18.15%  '
$ws.Range("D32").Value = '709.84'
$ws.Range("E32").Value = '  +5.56%  '
$ws.Range("D33").Value = '0.126'
$ws.Range("E33").Value = '  +13.36%  '
$ws.Range("D34").Value = '2.70'
$ws.Range("E34").Value = '  +5.50%  '
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").Value = '5.78'
$ws.Range("E36").Value = '  +42.58%  '
$ws.Range("D37").Value = '38.87'
$ws.Range("E37").Value = '  +5.45%  '
$ws.Range("D38").Value = '0.149'
$ws.Range("E38").Value = '  +0.08%  '
$ws.Range("D39").Value = '55.64'
$ws.Range("E39").Value = '  +3.09%  '
$ws.Range("D40").Value = '0.0470'
$ws.Range("E40").Value = '  +7.18%  '
$ws.Range("D41").Value = '0.0₃0731'
$ws.Range("E41").Value = '  +15.54%  '
$ws.Range("D42").Value = '2.90'
$ws.Range("E42").Value = '  +3.69%  '
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.32%  '
$ws.Range("D44").Value = '0.138'
$ws.Range("E44").Value = '  +4.53%  '
$ws.Range("D45").Value = '3.40'
$ws.Range("E45").Value = '  +9.78%  '
$ws.Range("D46").Value = '3.21'
$ws.Range("E46").Value = '  +3.71%  '
$ws.Range("D47").Value = '0.323'
$ws.Range("E47").Value = '  +16.47%  '
$ws.Range("D48").Value = '2.42'
$ws.Range("E48").Value = '  +42.71%  '
$ws.Range("D49").Value = '2.63'
$ws.Range("E49").Value = '  +7.08%  '
$ws.Range("D50").Value = '2.06'
$ws.Range("E50").Value = '  +5.41%  '
$ws.Range("D51").Value = '2.84'
$ws.Range("E51").Value = '  +3.86%  '
